# migrationStatus.xlsx update:
# "Added every element to db and created undoables"
#
# Column H lists TilerElements class files, column I tracks "migrated"
# status and column J tracks "undoable" status (x = done, * = pending).
# This change replaces the stale date placeholder in row 15 with the new
# "ConflictProfile.cs" element (marked pending in both status columns),
# fills in migrated/undoable status for every class that was previously
# missing it, and flips Reason.cs (row 40) from pending to done.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15 used to hold a stray date value; it now becomes a normal element
# row for the newly added ConflictProfile.cs, still pending migration.
$ws.Range("H15").Value = "ConflictProfile.cs"
$ws.Range("I15").Value = "*"
$ws.Range("J15").Value = "*"

# Elements newly marked fully migrated + undoable ("x") in both columns.
$xBothCols = @(13, 14, 38, 46, 47, 48)
foreach ($r in $xBothCols) {
    $ws.Cells.Item($r, 9).Value = "x"
    $ws.Cells.Item($r, 10).Value = "x"
}

# Elements newly marked pending ("*") in both columns.
$starBothCols = @(16, 17, 49, 50, 57, 58, 59, 61, 63, 65)
foreach ($r in $starBothCols) {
    $ws.Cells.Item($r, 9).Value = "*"
    $ws.Cells.Item($r, 10).Value = "*"
}

# Elements that already had a migrated ("I") status and now also get an
# undoable ("x") status in column J.
$xJOnly = @(18, 42, 43, 44, 45, 52, 53)
foreach ($r in $xJOnly) {
    $ws.Cells.Item($r, 10).Value = "x"
}

# Elements that already had a migrated ("I") status and now also get a
# pending ("*") undoable status in column J.
$starJOnly = @(51, 54, 55, 56, 62)
foreach ($r in $starJOnly) {
    $ws.Cells.Item($r, 10).Value = "*"
}

# Reason.cs (row 40) moves from pending ("*") to fully done ("x") in both
# the migrated and undoable columns.
$ws.Cells.Item(40, 9).Value = "x"
$ws.Cells.Item(40, 10).Value = "x"

# Restore the author's on-screen view/selection at the time of the commit.
$ws.Range("J15").Select() | Out-Null
